$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was collected. Insert a new row at position 18,
# which pushes the existing rows 18-27 down to 19-28 (dimension grows to R28).
$ws.Range("A18").EntireRow.Insert()

# Populate the newly inserted row 18 with the latest weekly data point.
$ws.Range("A18").Value = 11
$ws.Range("B18").Value = "Vega Monumental Concepción"
$ws.Range("C18").Value = "Bíobío"
$ws.Range("D18").Value = 44587
$ws.Range("E18").Value = 8
$ws.Range("F18").Value = 100112030
$ws.Range("G18").Value = "Poroto granado"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 220
$ws.Range("K18").Value = 23000
$ws.Range("L18").Value = 24000
$ws.Range("M18").Value = 23545
$ws.Range("N18").Value = "$/saco 25 kilos"
$ws.Range("O18").Value = "Región Metropolitana"
$ws.Range("P18").Value = 942
$ws.Range("Q18").Value = 25
$ws.Range("R18").Value = "Hortaliza"
